$wb = $excel.ActiveWorkbook

# The edit targets the "LoginWithInvalidCredentialsTest" worksheet (sheet1).
$ws = $wb.Worksheets.Item("LoginWithInvalidCredentialsTest")

# E2: 12345678955 -> 123456789
$ws.Range("E2").Value = 123456789

# E4: (empty) -> 123456789
$ws.Range("E4").Value = 123456789

# Update selection on that sheet to E4 (single cell) to match the saved view state,
# without disturbing which sheet is the active tab in the workbook.
$originalActive = $wb.ActiveSheet
$ws.Activate()
$ws.Range("E4").Select()
$originalActive.Activate()
